$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ("Row 57 Height: " + $ws.Rows(57).RowHeight)
Write-Host ("Cell E16 StyleIndex test")
$c = $ws.Range("E16")
Write-Host ("NumberFormat: " + $c.NumberFormat)
Write-Host ("Font Bold: " + $c.Font.Bold)
Write-Host ("Borders left: " + $c.Borders(7).LineStyle)

$c57 = $ws.Range("B57:J57")
Write-Host ("Row 57 address: " + $c57.Address())
